$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3703.7896
$ws.Range("I19").Value = 4844.4546
$ws.Range("J19").Value = 2135.375
$ws.Range("K19").Value = 4844.4546
$ws.Range("L19").Value = 2135.375
$ws.Range("M19").Value = -4669.4546
$ws.Range("N19").Value = -2485.375
$ws.Range("H53").Value = 67112.336
$ws.Range("I53").Value = 166930.17
$ws.Range("K53").Value = 166930.17
$ws.Range("M53").Value = -166293.17
$ws.Range("H112").Value = 1586.4
$ws.Range("J112").Value = 1863.2727
$ws.Range("L112").Value = 5589.8181
$ws.Range("N112").Value = -7805.8181
$ws.Range("H113").Value = 2908.0557
$ws.Range("I113").Value = 2012.375
$ws.Range("J113").Value = 3624.6
$ws.Range("K113").Value = 2012.375
$ws.Range("L113").Value = 3624.6
$ws.Range("M113").Value = 1241.625
$ws.Range("N113").Value = -10132.6
$ws.Range("H132").Value = 3402599.2
$ws.Range("I132").Value = 3969528.2
$ws.Range("K132").Value = 11908584.6
$ws.Range("M132").Value = -11906054.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H32").Value = 3762.845
$ws.Range("I32").Value = 3252.3823
$ws.Range("K32").Value = 3252.3823
$ws.Range("M32").Value = -2965.3823
$ws.Range("H45").Value = 20834800
$ws.Range("I45").Value = 41667988
$ws.Range("K45").Value = 41667988
$ws.Range("M45").Value = -41667611
$ws.Range("H61").Value = 2234.5715
$ws.Range("I61").Value = 1247.4546
$ws.Range("J61").Value = 3320.4
$ws.Range("K61").Value = 1247.4546
$ws.Range("L61").Value = 3320.4
$ws.Range("M61").Value = -1035.4546
$ws.Range("N61").Value = -3744.4
$ws.Range("H74").Value = 951
$ws.Range("I74").Value = 869.5789
$ws.Range("K74").Value = 869.5789
$ws.Range("M74").Value = 4.421100000000024
$ws.Range("H77").Value = 951
$ws.Range("I77").Value = 869.5789
$ws.Range("K77").Value = 4347.8945
$ws.Range("M77").Value = 20.10549999999967
$ws.Range("H132").Value = 6211.552
$ws.Range("I132").Value = 7256.95
$ws.Range("J132").Value = 3888.4443
$ws.Range("K132").Value = 21770.85
$ws.Range("L132").Value = 11665.3329
$ws.Range("M132").Value = -19240.85
$ws.Range("N132").Value = -16725.3329
$ws.Range("H136").Value = 2234.5715
$ws.Range("I136").Value = 1247.4546
$ws.Range("J136").Value = 3320.4
$ws.Range("K136").Value = 3742.3638
$ws.Range("L136").Value = 9961.200000000001
$ws.Range("M136").Value = -1192.3638
$ws.Range("N136").Value = -15061.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2132.1
$ws.Range("I99").Value = 2108.1667
$ws.Range("J99").Value = 2168
$ws.Range("K99").Value = 2108.1667
$ws.Range("L99").Value = 2168
$ws.Range("M99").Value = -610.1667000000002
$ws.Range("N99").Value = -5164
$ws.Range("H105").Value = 2832.1904
$ws.Range("I105").Value = 2804.7778
$ws.Range("K105").Value = 2804.7778
$ws.Range("M105").Value = -1057.7778
$ws.Range("H134").Value = 30183.172
$ws.Range("I134").Value = 38204.11
$ws.Range("J134").Value = 3112.5
$ws.Range("K134").Value = 114612.33
$ws.Range("L134").Value = 9337.5
$ws.Range("M134").Value = -112077.33
$ws.Range("N134").Value = -14407.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 413.83334
$ws.Range("I22").Value = 358.08334
$ws.Range("J22").Value = 525.3333
$ws.Range("K22").Value = 358.08334
$ws.Range("L22").Value = 525.3333
$ws.Range("M22").Value = -8.083340000000021
$ws.Range("N22").Value = -1225.3333
$ws.Range("H31").Value = 3548532.8
$ws.Range("I31").Value = 2300.2222
$ws.Range("J31").Value = 8335946.5
$ws.Range("K31").Value = 2300.2222
$ws.Range("L31").Value = 8335946.5
$ws.Range("M31").Value = -2005.2222
$ws.Range("N31").Value = -8336536.5
$ws.Range("H34").Value = 3548532.8
$ws.Range("I34").Value = 2300.2222
$ws.Range("J34").Value = 8335946.5
$ws.Range("K34").Value = 2300.2222
$ws.Range("L34").Value = 8335946.5
$ws.Range("M34").Value = -2098.2222
$ws.Range("N34").Value = -8336350.5
$ws.Range("H58").Value = 1413.0588
$ws.Range("I58").Value = 1415.8572
$ws.Range("J58").Value = 1400
$ws.Range("K58").Value = 1415.8572
$ws.Range("L58").Value = 1400
$ws.Range("M58").Value = -1212.8572
$ws.Range("N58").Value = -1806
$ws.Range("H99").Value = 2021.6061
$ws.Range("I99").Value = 1748.0454
$ws.Range("K99").Value = 1748.0454
$ws.Range("M99").Value = -250.0454
$ws.Range("H122").Value = 740.4483
$ws.Range("I122").Value = 680.625
$ws.Range("J122").Value = 1027.6
$ws.Range("K122").Value = 2041.875
$ws.Range("L122").Value = 3082.8
$ws.Range("M122").Value = 408.125
$ws.Range("N122").Value = -7982.799999999999
$ws.Range("H126").Value = 2021.6061
$ws.Range("I126").Value = 1748.0454
$ws.Range("K126").Value = 5244.1362
$ws.Range("M126").Value = -2774.1362
$ws.Range("H134").Value = 1091.5625
$ws.Range("I134").Value = 1079.6428
$ws.Range("K134").Value = 3238.9284
$ws.Range("M134").Value = -703.9284000000002
$ws.Range("H136").Value = 1413.0588
$ws.Range("I136").Value = 1415.8572
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 4247.571599999999
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -1697.571599999999
$ws.Range("N136").Value = -9300

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1012.13336
$ws.Range("J122").Value = 1420.2222
$ws.Range("L122").Value = 12781.9998
$ws.Range("N122").Value = -17681.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6632
$ws.Range("I80").Value = 6264
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 6264
$ws.Range("L80").Value = 7000
$ws.Range("M80").Value = -5266
$ws.Range("N80").Value = -8996
$ws.Range("H83").Value = 6632
$ws.Range("I83").Value = 6264
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 31320
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = -26328
$ws.Range("N83").Value = -44984
$ws.Range("H113").Value = 19231402
$ws.Range("I113").Value = 31250500
$ws.Range("K113").Value = 31250500
$ws.Range("M113").Value = -31248330
$ws.Range("H126").Value = 2814.7273
$ws.Range("I126").Value = 3602
$ws.Range("J126").Value = 1870
$ws.Range("K126").Value = 10806
$ws.Range("L126").Value = 5610
$ws.Range("M126").Value = -8336
$ws.Range("N126").Value = -10550
$ws.Range("H132").Value = 46032.26
$ws.Range("I132").Value = 78757.234
$ws.Range("J132").Value = 3489.8
$ws.Range("K132").Value = 236271.702
$ws.Range("L132").Value = 10469.4
$ws.Range("M132").Value = -233741.702
$ws.Range("N132").Value = -15529.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 537.7
$ws.Range("I22").Value = 595.5714
$ws.Range("J22").Value = 402.66666
$ws.Range("K22").Value = 595.5714
$ws.Range("L22").Value = 402.66666
$ws.Range("M22").Value = -300.5714
$ws.Range("N22").Value = -992.66666
$ws.Range("H27").Value = 537.7
$ws.Range("I27").Value = 595.5714
$ws.Range("J27").Value = 402.66666
$ws.Range("K27").Value = 595.5714
$ws.Range("L27").Value = 402.66666
$ws.Range("M27").Value = -488.5714
$ws.Range("N27").Value = -616.66666
$ws.Range("H55").Value = 233.93103
$ws.Range("I55").Value = 216.5
$ws.Range("J55").Value = 250.2
$ws.Range("K55").Value = 216.5
$ws.Range("L55").Value = 250.2
$ws.Range("M55").Value = -43.5
$ws.Range("N55").Value = -596.2
$ws.Range("H122").Value = 2817.8928
$ws.Range("I122").Value = 2804.7896
$ws.Range("J122").Value = 2845.5557
$ws.Range("K122").Value = 8414.3688
$ws.Range("L122").Value = 8536.667099999999
$ws.Range("M122").Value = -5964.3688
$ws.Range("N122").Value = -13436.6671
$ws.Range("H132").Value = 2427.2354
$ws.Range("I132").Value = 1495.25
$ws.Range("K132").Value = 4485.75
$ws.Range("M132").Value = -1955.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 419.58334
$ws.Range("I113").Value = 378.42105
$ws.Range("J113").Value = 576
$ws.Range("K113").Value = 1135.26315
$ws.Range("L113").Value = 1728
$ws.Range("M113").Value = 1034.73685
$ws.Range("N113").Value = -6068
$ws.Range("H114").Value = 29833.334
$ws.Range("J114").Value = 29833.334
$ws.Range("L114").Value = 29833.334
$ws.Range("N114").Value = -38511.334
$ws.Range("H132").Value = 2815.077
$ws.Range("I132").Value = 5152
$ws.Range("K132").Value = 15456
$ws.Range("M132").Value = -12926
$ws.Range("H136").Value = 4592.1514
$ws.Range("I136").Value = 5635.24
$ws.Range("K136").Value = 16905.72
$ws.Range("M136").Value = -14355.72
